# Generate Report for Handoff
#
# The localization-status workbook tracks, per source file, the
# handback/handoff status for each target language. The d036f29a...md
# file has just been re-handed-off for both zh-cn and de-de, so:
#   - its Status moves from "Handed back: in sync with en-US" to
#     "Ready for handoff" (reflected on the Overview sheet and on each
#     language sheet's Status column)
#   - each language sheet's "Latest Handoff Datetime" is stamped with
#     the new handoff time

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: Status columns for the zh-cn (B) and de-de (C) targets
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: Status + new handoff datetime for d036f29a...md (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("D3").Value = "2016-02-25 04:06:44"

# --- de-de sheet: Status + new handoff datetime for d036f29a...md (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("D3").Value = "2016-02-25 04:07:01"
